$wb = $excel.ActiveWorkbook

# sheet1 (展览)
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(4, 6).Value = 2117
$ws.Cells.Item(5, 6).Value = 9146
$ws.Cells.Item(7, 6).Value = 118
$ws.Cells.Item(10, 6).Value = 641
$ws.Cells.Item(11, 6).Value = 601
$ws.Cells.Item(14, 6).Value = 298
$ws.Cells.Item(16, 6).Value = 56
$ws.Cells.Item(17, 6).Value = 1512
$ws.Cells.Item(21, 6).Value = 1386
$ws.Cells.Item(23, 6).Value = 238
$ws.Cells.Item(25, 6).Value = 98
$ws.Cells.Item(27, 6).Value = 69
$ws.Cells.Item(28, 6).Value = 316
$ws.Cells.Item(29, 6).Value = 316
$ws.Cells.Item(30, 6).Value = 1074
$ws.Cells.Item(33, 6).Value = 234
$ws.Cells.Item(33, 7).Value = 29.9
$ws.Cells.Item(34, 6).Value = 205
$ws.Cells.Item(39, 6).Value = 137
$ws.Cells.Item(41, 6).Value = 158
$ws.Cells.Item(42, 6).Value = 70
$ws.Cells.Item(43, 6).Value = 503
$ws.Cells.Item(44, 6).Value = 1236
$ws.Cells.Item(45, 6).Value = 689
$ws.Cells.Item(46, 6).Value = 223

# sheet2 (演出)
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(6, 6).Value = 57
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(12, 6).Value = 227
$ws.Cells.Item(14, 6).Value = 162
$ws.Cells.Item(27, 6).Value = 236
$ws.Cells.Item(30, 6).Value = 232
$ws.Cells.Item(31, 6).Value = 5

# sheet3 (本地生活)
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(4, 6).Value = 752
$ws.Cells.Item(5, 6).Value = 320
$ws.Cells.Item(6, 6).Value = 147
$ws.Cells.Item(7, 6).Value = 2120
$ws.Cells.Item(8, 6).Value = 3188

# sheet4 (全部类型)
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(5, 6).Value = 752
$ws.Cells.Item(6, 6).Value = 9146
$ws.Cells.Item(7, 6).Value = 320
$ws.Cells.Item(8, 6).Value = 147
$ws.Cells.Item(10, 6).Value = 2120
$ws.Cells.Item(11, 6).Value = 3188
$ws.Cells.Item(12, 6).Value = 118
$ws.Cells.Item(14, 6).Value = 641
$ws.Cells.Item(16, 6).Value = 601
$ws.Cells.Item(18, 6).Value = 298
$ws.Cells.Item(19, 6).Value = 1512
$ws.Cells.Item(23, 6).Value = 1386
$ws.Cells.Item(24, 6).Value = 238
$ws.Cells.Item(25, 6).Value = 98
$ws.Cells.Item(27, 6).Value = 316
$ws.Cells.Item(28, 6).Value = 316
$ws.Cells.Item(29, 6).Value = 1074
$ws.Cells.Item(34, 6).Value = 234
$ws.Cells.Item(34, 7).Value = 29.9
$ws.Cells.Item(35, 6).Value = 205
$ws.Cells.Item(37, 6).Value = 236
$ws.Cells.Item(41, 6).Value = 137
$ws.Cells.Item(42, 6).Value = 232
$ws.Cells.Item(43, 6).Value = 158
$ws.Cells.Item(44, 6).Value = 70
$ws.Cells.Item(46, 6).Value = 503
$ws.Cells.Item(47, 6).Value = 689
